$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: reading ".Value" back from a cell in this host mis-resolves to the
# property descriptor text, not the cell's data - ".Value2" must be used
# whenever an existing cell's contents need to be fetched into a variable
# (or copied to another cell). Writes via ".Value =" are unaffected.

# --- Header row 3: insert "Amount" in D, shift old "Price" header to E ---
$ws.Range("E3").Value = $ws.Range("D3").Value2
$ws.Range("D3").Value = "Amount"

# --- Row 4 (PISO-64-2-C-7-2-FB): was a text "260 + GST" in D4 ---
# Split into Amount (D4) / Price (E4) and add the Amount*Price formula in F4.
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 260
$ws.Range("F4").Formula = "=D4*E4"

# --- Row 5 (KMSS_M): old D5 held the Price; move it to E5, set Amount=1 ---
$ws.Range("E5").Value = $ws.Range("D5").Value2
$ws.Range("D5").Value = 1

# --- Row 6 (POLARIS-K05) ---
$ws.Range("E6").Value = $ws.Range("D6").Value2
$ws.Range("D6").Value = 1

# --- Row 7 (C240TME-1064) ---
$ws.Range("E7").Value = $ws.Range("D7").Value2
$ws.Range("D7").Value = 1

# --- Row 8 (WW11050-C14) ---
$ws.Range("E8").Value = $ws.Range("D8").Value2
$ws.Range("D8").Value = 1

# --- Row 9 (GR13-1210): E9 held a stray "41 deg" label -> replace with Price ---
$ws.Range("E9").Value = $ws.Range("D9").Value2
$ws.Range("D9").Value = 1

# --- Row 10 (1740-1053-012): E10 held a stray "70 deg" label -> replace with Price ---
$ws.Range("E10").Value = $ws.Range("D10").Value2
$ws.Range("D10").Value = 1

# --- Row 11 (PA4FK): 4 units at 31 each ---
$ws.Range("E11").Value = $ws.Range("D11").Value2
$ws.Range("D11").Value = 4

# --- Add the shared Amount*Price formula across F5:F17 (rows 12-17 are ---
# --- currently blank Amount/Price cells, so those formulas evaluate to 0) ---
$ws.Range("F5:F17").Formula = "=D5*E5"
$ws.Range("F5:F17").Style = "Normal"

# Row 12's formula wasn't grouped into the shared block (matches the source
# workbook, where F12 carries its own literal formula).
$ws.Range("F12").Formula = "=D12*E12"
$ws.Range("F12").Style = "Normal"

# --- Row 19: move "Total" from C19 to E19, clear the old SUM, add new one ---
$ws.Range("E19").Value = "Total"
$ws.Range("E19").Font.Bold = $true
$ws.Range("C19").Clear()
$ws.Range("D19").Clear()
$ws.Range("F19").Formula = "=SUM(F4:F18)"

# --- Selection moved from B18 to B17 ---
$ws.Range("B17").Select()
